$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @("hexadecanoic acid", "hexadecanoic acid", "C16H32O2", "CCCCCCCCCCCCCCCC(=O)O", 256.42, 6.4, 16, 0, 32, 2, 0.7494579205990172, 0, 0.125793619842446, 0.1247874580765931, 15, 0, 0, 0, 1, 0, 0, 0, 0.8244793697839481, 0, 0, 0, 0.1755596287341081, 0, 0, 0, 1.000015722652411)
$row3 = @("2,5-hexanedione", "hexane-2,5-dione", "C6H10O2", "CC(=O)CCC(=O)C", 114.14, -0.3, 6, 0, 10, 2, 0.6313825127036973, 0, 0.08831259856316805, 0.2803399334151043, 0, 0, 0, 0, 0, 0, 0, 2, 0, 0, 0, 0, 0, 0, 0, 1.000035044681969, 1.000015722652411)
$row4 = @("trans-2-pentenoic acid", "(e)-pent-2-enoic acid", "C5H8O2", "CCC=CC(=O)O", 100.12, 1, 5, 0, 8, 2, 0.5998302037554933, 0, 0.08054334798242109, 0.3195964842189373, 4, 0, 0, 0, 1, 0, 0, 0, 0.5503395924890131, 0, 0, 0, 0.4496304434678386, 0, 0, 0, 1.000015722652411)
$row5 = @("9,12-octadecadienoic acid (z,z)-", "(9z,12z)-octadeca-9,12-dienoic acid", "C18H32O2", "CCCCCC=CCC=CCCCCCCCC(=O)O", 280.4, 6.8, 18, 0, 32, 2, 0.7710342368045648, 0, 0.1150356633380885, 0.1141155492154066, 17, 0, 0, 0, 1, 0, 0, 0, 0.8396398002853066, 0, 0, 0, 0.1605456490727532, 0, 0, 0, 1.000015722652411)
$row6 = @("2-cyclopenten-1-one, 2-methyl-", "2-methylcyclopent-2-en-1-one", "C6H8O", "CC1=CCCC1=O", 96.13, 0.9, 6, 0, 8, 1, 0.7496723187350464, 0, 0.08388640382814938, 0.1664308748569645, 3, 0, 0, 0, 0, 0, 0, 1, 0.4377509622386352, 0, 0, 0, 0, 0, 0, 0.5936960366170811, 1.000015722652411)
$row7 = @("1-hexene, 4,5-dimethyl-", "4,5-dimethylhex-1-ene", "C8H16", "CC(C)C(C)CC=C", 112.21, 3.5, 8, 0, 16, 0, 0.8563229658675697, 0, 0.1437305053025577, 0, 8, 0, 0, 0, 0, 0, 0, 0, 1.000053471170127, 0, 0, 0, 0, 0, 0, 0, 1.000015722652411)
$row8 = @("n-hexadecanoic acid", "hexadecanoic acid", "C16H32O2", "CCCCCCCCCCCCCCCC(=O)O", 256.42, 6.4, 16, 0, 32, 2, 0.7494579205990172, 0, 0.125793619842446, 0.1247874580765931, 15, 0, 0, 0, 1, 0, 0, 0, 0.8244793697839481, 0, 0, 0, 0.1755596287341081, 0, 0, 0, 1.000015722652411)
$row9 = @("2,4,5-trichlorophenol", "2,4,5-trichlorophenol", "C6H3Cl3O", "C1=C(C(=CC(=C1Cl)Cl)Cl)O", 197.4, 3.7, 6, 3, 3, 1, 0.3650759878419453, 0.5387537993920973, 0.01531914893617021, 0.08104863221884498, 0, 6, 3, 1, 0, 0, 3, 0, 0, 0.3752887537993921, 0.5387537993920973, 0.08615501519756839, 0, 0, 0.5387537993920973, 0, 1.000015722652411)
$row10 = @("(z)-octadec-9-enoic acid", "(z)-octadec-9-enoic acid", "C18H34O2", "CCCCCCCCC=CCCCCCCCC(=O)O", 282.5, 6.5, 18, 0, 34, 2, 0.7653026548672566, 0, 0.121316814159292, 0.1132672566371681, 17, 0, 0, 0, 1, 0, 0, 0, 0.8405345132743363, 0, 0, 0, 0.1593522123893805, 0, 0, 0, 1.000015722652411)
$row11 = @("tetradecanoic acid", "tetradecanoic acid", "C14H28O2", "CCCCCCCCCCCCCC(=O)O", 228.37, 5.3, 14, 0, 28, 2, 0.7363226343214958, 0, 0.1235889127293427, 0.1401147261023777, 13, 0, 0, 0, 1, 0, 0, 0, 0.8029031834303979, 0, 0, 0, 0.1971230897228183, 0, 0, 0, 1.000015722652411)
$row12 = @("2-butanone", "butan-2-one", "C4H8O", "CCC(=O)C", 72.11, 0.3, 4, 0, 8, 1, 0.6662598807377618, 0, 0.1118291499098599, 0.2218693662460131, 1, 0, 0, 0, 0, 0, 0, 1, 0.2085009014006379, 0, 0, 0, 0, 0, 0, 0.7914574954929968, 1.000015722652411)
$row13 = @("oleic acid", "(z)-octadec-9-enoic acid", "C18H34O2", "CCCCCCCCC=CCCCCCCCC(=O)O", 282.5, 6.5, 18, 0, 34, 2, 0.7653026548672566, 0, 0.121316814159292, 0.1132672566371681, 17, 0, 0, 0, 1, 0, 0, 0, 0.8405345132743363, 0, 0, 0, 0.1593522123893805, 0, 0, 0, 1.000015722652411)
$row14 = @("n-decanoic acid", "decanoic acid", "C10H20O2", "CCCCCCCCCC(=O)O", 172.26, 4.1, 10, 0, 20, 2, 0.6972599558806455, 0, 0.1170323928944619, 0.1857540926506444, 9, 0, 0, 0, 1, 0, 0, 0, 0.7387147335423198, 0, 0, 0, 0.2613317078834321, 0, 0, 0, 1.000015722652411)
$row15 = @("2-methylcyclopent-2-en-1-one", "2-methylcyclopent-2-en-1-one", "C6H8O", "CC1=CCCC1=O", 96.13, 0.9, 6, 0, 8, 1, 0.7496723187350464, 0, 0.08388640382814938, 0.1664308748569645, 3, 0, 0, 0, 0, 0, 0, 1, 0.4377509622386352, 0, 0, 0, 0, 0, 0, 0.5936960366170811, 1.000015722652411)
$row16 = @("phenol", "phenol", "C6H6O", "C1=CC=C(C=C1)O", 94.11, 1.5, 6, 0, 6, 1, 0.765763468281798, 0, 0.06426522154925088, 0.1700031877590054, 0, 6, 0, 1, 0, 0, 0, 0, 0, 0.8193178195728402, 0, 0.1807140580172139, 0, 0, 0, 0, 1.000015722652411)
$row17 = @("(9z,12z)-octadeca-9,12-dienoic acid", "(9z,12z)-octadeca-9,12-dienoic acid", "C18H32O2", "CCCCCC=CCC=CCCCCCCCC(=O)O", 280.4, 6.8, 18, 0, 32, 2, 0.7710342368045648, 0, 0.1150356633380885, 0.1141155492154066, 17, 0, 0, 0, 1, 0, 0, 0, 0.8396398002853066, 0, 0, 0, 0.1605456490727532, 0, 0, 0, 1.000015722652411)
$row18 = @("oxacycloheptadecan-2-one", "oxacycloheptadecan-2-one", "C16H30O2", "C1CCCCCCCC(=O)OCCCCCCC1", 254.41, 6.3, 16, 0, 30, 2, 0.7553791124562713, 0, 0.1188632522306513, 0.1257733579654888, 14, 0, 0, 0, 0, 1, 0, 0, 0.771895758814512, 0, 0, 0, 0, 0.2281199638378994, 0, 0, 1.000015722652411)

$rows = $row2, $row3, $row4, $row5, $row6, $row7, $row8, $row9, $row10, $row11, $row12, $row13, $row14, $row15, $row16, $row17, $row18

for ($i = 0; $i -lt $rows.Length; $i++) {
    $row = $rows[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($i + 2, $j + 1).Value = $row[$j]
    }
}

Write-Host "done writing rows"